$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4397.5
$ws.Range("I62").Value = 2082.1428
$ws.Range("J62").Value = 4860.5713
$ws.Range("K62").Value = 2082.1428
$ws.Range("L62").Value = 4860.5713
$ws.Range("M62").Value = -1458.1428
$ws.Range("N62").Value = -6108.5713
$ws.Range("H65").Value = 4397.5
$ws.Range("I65").Value = 2082.1428
$ws.Range("J65").Value = 4860.5713
$ws.Range("K65").Value = 10410.714
$ws.Range("L65").Value = 24302.8565
$ws.Range("M65").Value = -7290.714
$ws.Range("N65").Value = -30542.8565
$ws.Range("H101").Value = 1743.4
$ws.Range("I101").Value = 224
$ws.Range("J101").Value = 4022.5
$ws.Range("K101").Value = 672
$ws.Range("L101").Value = 12067.5
$ws.Range("M101").Value = 950
$ws.Range("N101").Value = -15311.5
$ws.Range("H137").Value = 20116016
$ws.Range("I137").Value = 3704650.8
$ws.Range("J137").Value = 76924590
$ws.Range("K137").Value = 11113952.4
$ws.Range("L137").Value = 230773770
$ws.Range("M137").Value = -11111402.4
$ws.Range("N137").Value = -230778870
$ws.Range("H138").Value = 2230.3052
$ws.Range("I138").Value = 1727.5714
$ws.Range("J138").Value = 2963.4583
$ws.Range("K138").Value = 5182.7142
$ws.Range("L138").Value = 8890.374899999999
$ws.Range("M138").Value = -42.71420000000035
$ws.Range("N138").Value = -19170.3749
$ws.Range("H141").Value = 1166.2084
$ws.Range("I141").Value = 726.7895
$ws.Range("J141").Value = 2836
$ws.Range("K141").Value = 2180.3685
$ws.Range("L141").Value = 8508
$ws.Range("M141").Value = 2999.6315
$ws.Range("N141").Value = -18868

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3429.5361
$ws.Range("I32").Value = 3507.1462
$ws.Range("J32").Value = 3005.2666
$ws.Range("K32").Value = 3507.1462
$ws.Range("L32").Value = 3005.2666
$ws.Range("M32").Value = -3220.1462
$ws.Range("N32").Value = -3579.2666
$ws.Range("H45").Value = 46405.09
$ws.Range("I45").Value = 200456.4
$ws.Range("J45").Value = 1095.8823
$ws.Range("K45").Value = 200456.4
$ws.Range("L45").Value = 1095.8823
$ws.Range("M45").Value = -200079.4
$ws.Range("N45").Value = -1849.8823
$ws.Range("H61").Value = 17780460
$ws.Range("I61").Value = 18521288
$ws.Range("J61").Value = 600
$ws.Range("K61").Value = 18521288
$ws.Range("L61").Value = 600
$ws.Range("M61").Value = -18521076
$ws.Range("N61").Value = -1024
$ws.Range("H74").Value = 19500302
$ws.Range("I74").Value = 27193594
$ws.Range("J74").Value = 10628.2
$ws.Range("K74").Value = 27193594
$ws.Range("L74").Value = 10628.2
$ws.Range("M74").Value = -27192720
$ws.Range("N74").Value = -12376.2
$ws.Range("H77").Value = 19500302
$ws.Range("I77").Value = 27193594
$ws.Range("J77").Value = 10628.2
$ws.Range("K77").Value = 135967970
$ws.Range("L77").Value = 53141
$ws.Range("M77").Value = -135963602
$ws.Range("N77").Value = -61877
$ws.Range("H132").Value = 827014.6
$ws.Range("I132").Value = 1114270.1
$ws.Range("J132").Value = 84938
$ws.Range("K132").Value = 3342810.3
$ws.Range("L132").Value = 254814
$ws.Range("M132").Value = -3340280.3
$ws.Range("N132").Value = -259874
$ws.Range("H136").Value = 17780460
$ws.Range("I136").Value = 18521288
$ws.Range("J136").Value = 600
$ws.Range("K136").Value = 55563864
$ws.Range("L136").Value = 1800
$ws.Range("M136").Value = -55561314
$ws.Range("N136").Value = -6900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1953.1754
$ws.Range("I86").Value = 1747.8948
$ws.Range("J86").Value = 2363.7368
$ws.Range("K86").Value = 1747.8948
$ws.Range("L86").Value = 2363.7368
$ws.Range("M86").Value = -624.8948
$ws.Range("N86").Value = -4609.736800000001
$ws.Range("H89").Value = 1953.1754
$ws.Range("I89").Value = 1747.8948
$ws.Range("J89").Value = 2363.7368
$ws.Range("K89").Value = 8739.474
$ws.Range("L89").Value = 11818.684
$ws.Range("M89").Value = -3123.474
$ws.Range("N89").Value = -23050.684
$ws.Range("H94").Value = 860.8333
$ws.Range("I94").Value = 623.3333
$ws.Range("J94").Value = 1098.3334
$ws.Range("K94").Value = 623.3333
$ws.Range("L94").Value = 1098.3334
$ws.Range("M94").Value = -172.3333
$ws.Range("N94").Value = -2000.3334
$ws.Range("H105").Value = 1751.1765
$ws.Range("I105").Value = 1891.1111
$ws.Range("J105").Value = 1593.75
$ws.Range("K105").Value = 1891.1111
$ws.Range("L105").Value = 1593.75
$ws.Range("M105").Value = -144.1111000000001
$ws.Range("N105").Value = -5087.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1602.2823
$ws.Range("I58").Value = 729.13336
$ws.Range("J58").Value = 3697.84
$ws.Range("K58").Value = 729.13336
$ws.Range("L58").Value = 3697.84
$ws.Range("M58").Value = -526.13336
$ws.Range("N58").Value = -4103.84
$ws.Range("H136").Value = 1602.2823
$ws.Range("I136").Value = 729.13336
$ws.Range("J136").Value = 3697.84
$ws.Range("K136").Value = 2187.40008
$ws.Range("L136").Value = 11093.52
$ws.Range("M136").Value = 362.5999199999997
$ws.Range("N136").Value = -16193.52

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 113599.4
$ws.Range("J141").Value = 113599.4
$ws.Range("L141").Value = 113599.4
$ws.Range("N141").Value = -123959.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1781.5834
$ws.Range("I61").Value = 1719
$ws.Range("J61").Value = 1844.1666
$ws.Range("K61").Value = 1719
$ws.Range("L61").Value = 1844.1666
$ws.Range("M61").Value = -1517
$ws.Range("N61").Value = -2248.1666
$ws.Range("H113").Value = 1781.5834
$ws.Range("I113").Value = 1719
$ws.Range("J113").Value = 1844.1666
$ws.Range("K113").Value = 1719
$ws.Range("L113").Value = 1844.1666
$ws.Range("M113").Value = 451
$ws.Range("N113").Value = -6184.1666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4138.8335
$ws.Range("I132").Value = 4419.2563
$ws.Range("J132").Value = 493.33334
$ws.Range("K132").Value = 13257.7689
$ws.Range("L132").Value = 1480.00002
$ws.Range("M132").Value = -10727.7689
$ws.Range("N132").Value = -6540.000019999999
$ws.Range("H136").Value = 5840.263
$ws.Range("I136").Value = 6589.242
$ws.Range("J136").Value = 897
$ws.Range("K136").Value = 19767.726
$ws.Range("L136").Value = 2691
$ws.Range("M136").Value = -17217.726
$ws.Range("N136").Value = -7791
